$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must remain TEXT (not auto-converted
# to numbers by Excel). Force the whole D range to Text format before writing, then restore
# the cell style back to Normal/General so no stray per-cell style survives in the saved file.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.864.27"
$ws.Range("D3").Value = "1.871.42"
$ws.Range("D5").Value = "0.7349"
$ws.Range("D6").Value = "241.94"
$ws.Range("D7").Value = "0.9996"
$ws.Range("D8").Value = "0.3158"
$ws.Range("D9").Value = "24.65"
$ws.Range("D10").Value = "0.07097"
$ws.Range("D11").Value = "0.08471"
$ws.Range("D12").Value = "0.7506"
$ws.Range("D13").Value = "5.375"
$ws.Range("D14").Value = "1.864.17"
$ws.Range("D15").Value = "92.43"
$ws.Range("D16").Value = "29.858.89"
$ws.Range("D17").Value = "6.042"
$ws.Range("D18").Value = "13.57"
$ws.Range("D19").Value = "243.19"
$ws.Range("D20").Value = "0.000007824"
$ws.Range("D21").Value = "0.9992"
$ws.Range("D22").Value = "2.124.37"
$ws.Range("D23").Value = "7.918"
$ws.Range("D24").Value = "1.000"
$ws.Range("D25").Value = "0.1565"
$ws.Range("D26").Value = "9.321"
$ws.Range("D27").Value = "164.18"
$ws.Range("D28").Value = "18.63"
$ws.Range("D29").Value = "2.024"
$ws.Range("D30").Value = "1.464"
$ws.Range("D31").Value = "4.539"
$ws.Range("D32").Value = "1.530"
$ws.Range("D33").Value = "4.267"
$ws.Range("D34").Value = "0.05325"
$ws.Range("D35").Value = "1.234"
$ws.Range("D36").Value = "0.7491"
$ws.Range("D37").Value = "0.9980"
$ws.Range("D38").Value = "2.690"
$ws.Range("D39").Value = "0.01948"
$ws.Range("D40").Value = "2.757"
$ws.Range("D41").Value = "0.4465"
$ws.Range("D42").Value = "1.097.73"
$ws.Range("D43").Value = "6.052"
$ws.Range("D44").Value = "72.25"
$ws.Range("D45").Value = "0.8674"
$ws.Range("D46").Value = "1.001"
$ws.Range("D47").Value = "7.720"
$ws.Range("D48").Value = "102.56"
$ws.Range("D49").Value = "3.068"
$ws.Range("D50").Value = "1.839"
$ws.Range("D51").Value = "2.020.27"

$dRange.Style = "Normal"

# Column E holds percentage-change strings (never ambiguous as a pure number because of the
# leading/trailing double-spaces and "%" sign), so they can be written directly as text.
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  -4.42%  "
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  +4.76%  "
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("E51").Value = "  +0.37%  "
